$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B-column text values where changed
$ws.Range("B5").Value = "<an>"
$ws.Range("B7").Value = "<light>"
$ws.Range("B9").Value = "<bow>"
$ws.Range("B13").Value = "<then>"

# Update C-column numeric values for rows 2-18
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 11
$ws.Range("C4").Value = 9
$ws.Range("C5").Value = 12
$ws.Range("C6").Value = 14
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 12
$ws.Range("C9").Value = 17
$ws.Range("C10").Value = 11
$ws.Range("C11").Value = 8
$ws.Range("C12").Value = 11
$ws.Range("C13").Value = 13
$ws.Range("C14").Value = 11
$ws.Range("C15").Value = 10
$ws.Range("C16").Value = 9
$ws.Range("C17").Value = 14
$ws.Range("C18").Value = 11
